# Auto-generated edit script: updates computed price/profit columns (H-N)
# across several sheets, per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -526

$ws.Range("H18").Value = 315.5
$ws.Range("I18").Value = 332
$ws.Range("J18").Value = 299
$ws.Range("K18").Value = 332
$ws.Range("L18").Value = 299
$ws.Range("M18").Value = -48
$ws.Range("N18").Value = -867

$ws.Range("H111").Value = 1253.6666
$ws.Range("I111").Value = 229
$ws.Range("J111").Value = 1766
$ws.Range("K111").Value = 687
$ws.Range("L111").Value = 5298
$ws.Range("M111").Value = 2380
$ws.Range("N111").Value = -11432

$ws.Range("H121").Value = 666.6667
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 666.6667
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2000.0001
$ws.Range("N121").Value = -5494.0001

$ws.Range("H137").Value = 71429650
$ws.Range("I137").Value = 76924170
$ws.Range("J137").Value = 900
$ws.Range("K137").Value = 230772510
$ws.Range("L137").Value = 2700
$ws.Range("M137").Value = -230769960
$ws.Range("N137").Value = -7800

$ws.Range("H141").Value = 3903.6365
$ws.Range("I141").Value = 4573.75
$ws.Range("J141").Value = 3520.7144
$ws.Range("K141").Value = 13721.25
$ws.Range("L141").Value = 10562.1432
$ws.Range("M141").Value = -8541.25
$ws.Range("N141").Value = -20922.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10022.4
$ws.Range("I74").Value = 2276
$ws.Range("J74").Value = 31325
$ws.Range("K74").Value = 2276
$ws.Range("L74").Value = 31325
$ws.Range("M74").Value = -1402
$ws.Range("N74").Value = -33073

$ws.Range("H77").Value = 10022.4
$ws.Range("I77").Value = 2276
$ws.Range("J77").Value = 31325
$ws.Range("K77").Value = 11380
$ws.Range("L77").Value = 156625
$ws.Range("M77").Value = -7012
$ws.Range("N77").Value = -165361

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 565.125
$ws.Range("I64").Value = 669.6667
$ws.Range("J64").Value = 502.4
$ws.Range("K64").Value = 669.6667
$ws.Range("L64").Value = 502.4
$ws.Range("M64").Value = -444.6667
$ws.Range("N64").Value = -952.4

$ws.Range("H67").Value = 565.125
$ws.Range("I67").Value = 669.6667
$ws.Range("J67").Value = 502.4
$ws.Range("K67").Value = 669.6667
$ws.Range("L67").Value = 502.4
$ws.Range("M67").Value = 110.3333
$ws.Range("N67").Value = -2062.4

$ws.Range("H86").Value = 9249.857
$ws.Range("I86").Value = 1742.6666
$ws.Range("J86").Value = 22762.8
$ws.Range("K86").Value = 1742.6666
$ws.Range("L86").Value = 22762.8
$ws.Range("M86").Value = -619.6666
$ws.Range("N86").Value = -25008.8

$ws.Range("H89").Value = 9249.857
$ws.Range("I89").Value = 1742.6666
$ws.Range("J89").Value = 22762.8
$ws.Range("K89").Value = 8713.333000000001
$ws.Range("L89").Value = 113814
$ws.Range("M89").Value = -3097.333000000001
$ws.Range("N89").Value = -125046

$ws.Range("H94").Value = 633.3333
$ws.Range("I94").Value = 520
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 520
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -69
$ws.Range("N94").Value = -2102

$ws.Range("H107").Value = 461.2353
$ws.Range("I107").Value = 247.71428
$ws.Range("J107").Value = 610.7
$ws.Range("K107").Value = 247.71428
$ws.Range("L107").Value = 610.7
$ws.Range("M107").Value = 1672.28572
$ws.Range("N107").Value = -4450.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1137.7
$ws.Range("I31").Value = 1076.1072
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1076.1072
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -781.1071999999999
$ws.Range("N31").Value = -2590

$ws.Range("H34").Value = 1137.7
$ws.Range("I34").Value = 1076.1072
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1076.1072
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -874.1071999999999
$ws.Range("N34").Value = -2404

$ws.Range("H41").Value = 7998.3335
$ws.Range("I41").Value = 1997.5
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 1997.5
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -1569.5
$ws.Range("N41").Value = -20856

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 335.5
$ws.Range("I11").Value = 335.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1006.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -866.5
$ws.Range("N11").ClearContents()

$ws.Range("H121").Value = 687.0909
$ws.Range("I121").Value = 244.28572
$ws.Range("J121").Value = 893.73334
$ws.Range("K121").Value = 732.85716
$ws.Range("L121").Value = 2681.20002
$ws.Range("M121").Value = 577.14284
$ws.Range("N121").Value = -5301.20002

$ws.Range("H140").Value = 6328.364
$ws.Range("I140").Value = 12030
$ws.Range("J140").Value = 2381.077
$ws.Range("K140").Value = 36090
$ws.Range("L140").Value = 7143.231000000001
$ws.Range("M140").Value = -30910
$ws.Range("N140").Value = -17503.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 186.89473
$ws.Range("I107").Value = 219.25
$ws.Range("J107").Value = 163.36363
$ws.Range("K107").Value = 219.25
$ws.Range("L107").Value = 163.36363
$ws.Range("M107").Value = 1700.75
$ws.Range("N107").Value = -4003.36363

$ws.Range("H131").Value = 35000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 35000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 366.82608
$ws.Range("I46").Value = 370.76923
$ws.Range("J46").Value = 361.7
$ws.Range("K46").Value = 370.76923
$ws.Range("L46").Value = 361.7
$ws.Range("M46").Value = -182.76923
$ws.Range("N46").Value = -737.7

$ws.Range("H68").Value = 1372.8182
$ws.Range("I68").Value = 1360.1
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 1360.1
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -611.0999999999999
$ws.Range("N68").Value = -2998

$ws.Range("H71").Value = 1372.8182
$ws.Range("I71").Value = 1360.1
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 6800.5
$ws.Range("L71").Value = 7500
$ws.Range("M71").Value = -3056.5
$ws.Range("N71").Value = -14988

$ws.Range("H93").Value = 585.4828
$ws.Range("I93").Value = 562.125
$ws.Range("J93").Value = 697.6
$ws.Range("K93").Value = 562.125
$ws.Range("L93").Value = 697.6
$ws.Range("M93").Value = 685.875
$ws.Range("N93").Value = -3193.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14723.637
$ws.Range("I62").Value = 17509
$ws.Range("J62").Value = 7296
$ws.Range("K62").Value = 17509
$ws.Range("L62").Value = 7296
$ws.Range("M62").Value = -16885
$ws.Range("N62").Value = -8544

$ws.Range("H65").Value = 14723.637
$ws.Range("I65").Value = 17509
$ws.Range("J65").Value = 7296
$ws.Range("K65").Value = 87545
$ws.Range("L65").Value = 36480
$ws.Range("M65").Value = -84425
$ws.Range("N65").Value = -42720

$ws.Range("H81").Value = 4700
$ws.Range("I81").Value = 1700
$ws.Range("J81").Value = 4900
$ws.Range("K81").Value = 3400
$ws.Range("L81").Value = 9800
$ws.Range("M81").Value = -2339
$ws.Range("N81").Value = -11922

$ws.Range("H84").Value = 4700
$ws.Range("I84").Value = 1700
$ws.Range("J84").Value = 4900
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 49000
$ws.Range("M84").Value = -11696
$ws.Range("N84").Value = -59608

